# Remove the "Install pushover" block (pushover library install removed).
#
# The document has a list of "Install <thing>" steps. The block being
# removed looks like:
#
#   Install pushover
#   cd ~
#   git clone https://github.com/thibauth/python-pushover.git
#   cd python-pushover
#   sudo pip3 install .
#   <blank paragraph>
#
# Right after it comes the "Install smbus" block. We delete the pushover
# block in its entirety (including the trailing blank separator paragraph,
# since the block before it already ends with its own blank separator), so
# the document flows directly from "sudo apt install git" / <blank> into
# "Install smbus" just like every other step in this list.
#
# We locate the block by searching for the literal command text rather
# than hard-coded paragraph indices, so the script is resilient to the
# exact paragraph numbering.

$d = $word.ActiveDocument

$startRange = $d.Content
$startRange.Find.ClearFormatting()
$found = $startRange.Find.Execute("Install pushover", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$endRange = $d.Content
$endRange.Find.ClearFormatting()
$found2 = $endRange.Find.Execute("sudo pip3 install .", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Expand the end of the deletion through the paragraph mark that follows
# "sudo pip3 install ." (the blank separator paragraph), so exactly one
# blank line remains between "sudo apt install git" and "Install smbus".
$blockStart = $startRange.Paragraphs(1).Range.Start
$cmdParagraph = $endRange.Paragraphs(1)
$blankParagraph = $cmdParagraph.Next()
$blockEnd = $blankParagraph.Range.End

$deleteRange = $d.Range($blockStart, $blockEnd)
$deleteRange.Delete()

# The stray "_GoBack" bookmark that Word leaves behind from the last edit
# location needs to move to the new last-edit point: right at the very
# start of the paragraph that now begins the next block ("Install smbus").
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$smbusRange = $d.Content
$smbusRange.Find.ClearFormatting()
$smbusRange.Find.Execute("Install smbus", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$smbusParagraph = $smbusRange.Paragraphs(1)
$goBackPoint = $d.Range($smbusParagraph.Range.Start, $smbusParagraph.Range.Start)
$d.Bookmarks.Add("_GoBack", $goBackPoint)

Write-Output "paragraphs=$($d.Paragraphs.Count)"
Write-Output "bookmark exists=$($d.Bookmarks.Exists('_GoBack'))"
